$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the two broken/renamed SUM formulas on row 18 (French "SOMME" instead
#    of "SUM" -> still broken/#NAME?, matching the author's actual edit)
# ---------------------------------------------------------------------------
$ws.Range("I18").Formula = "=SOMME(I7:I17)"
$ws.Range("J18").Formula = "=SOMME(J7:J17)"

# ---------------------------------------------------------------------------
# 2. Extend the "DEPENSES AUTRES" block (rows 21-25) with the same formatting
#    already used by the "PAROISSIENS" block (rows 7-17): bordered cells,
#    centered bold category codes, bordered amount columns.
# ---------------------------------------------------------------------------
$ws.Range("F7:G7").Copy()
$ws.Range("F21:G25").PasteSpecial(-4122)

$ws.Range("I7:J7").Copy()
$ws.Range("I21:J25").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Apply a dollar amount number format to both amount blocks.
# ---------------------------------------------------------------------------
$ws.Range("I7:J17").NumberFormat = '#,##0.00\$'
$ws.Range("I21:J25").NumberFormat = '#,##0.00\$'

# ---------------------------------------------------------------------------
# 4. Row 26: "TOTAL REVENUS DES AUTRES" with totals for the second block,
#    using the same style already used for the "TOTAL REVENUS DES
#    PAROISSIENS" caption (row 18/20).
# ---------------------------------------------------------------------------
$ws.Range("F20").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("F26").Value = "TOTAL REVENUS DES AUTRES"
$ws.Rows.Item(26).RowHeight = 15.75

$ws.Range("I26").Formula = "=SOMME(I21:I25)"
$ws.Range("J26").Formula = "=SOMME(J21:J25)"

# ---------------------------------------------------------------------------
# 5. Row 28: "GRAND  TOTAL REVENUS" big bold title + grand totals.
# ---------------------------------------------------------------------------
$ws.Range("F20").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F28").Font.Size = 16
$ws.Range("F28").Value = "GRAND  TOTAL REVENUS "
$ws.Rows.Item(28).RowHeight = 21

$ws.Range("I28").Formula = "=SOMME(I21:I25)"
$ws.Range("J28").Formula = "=SOMME(J21:J25)"

# ---------------------------------------------------------------------------
# 6. Selection matches where the author ended up after the edit.
# ---------------------------------------------------------------------------
$ws.Range("J27").Select()

$wb.Save()
